$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.270956
$ws.Range("H2").Value = 0.812868
$ws.Range("I2").Value = 0.05752762181187367
$ws.Range("J2").Value = 0.05752762181187366
$ws.Range("M2").Value = 0.8077876666666667
$ws.Range("N2").Value = 2.423363
$ws.Range("O2").Value = 0.04902038147436601
$ws.Range("P2").Value = 0.04902038147436601
$ws.Range("Q2").Value = 0.2188749150093334
$ws.Range("R2").Value = 1.969874235084
$ws.Range("S2").Value = 0.002820025966531106
$ws.Range("T2").Value = 0.002820025966531106
$ws.Range("G3").Value = 0.270956
$ws.Range("H3").Value = 0.812868
$ws.Range("I3").Value = 0.05752762181187367
$ws.Range("J3").Value = 0.05752762181187366
$ws.Range("O3").Value = 0.4722854529078861
$ws.Range("P3").Value = 0.4722854529078861
$ws.Range("Q3").Value = 2.108744062292
$ws.Range("R3").Value = 18.978696560628
$ws.Range("S3").Value = 0.02716945892213434
$ws.Range("T3").Value = 0.02716945892213434
$ws.Range("G4").Value = 0.270956
$ws.Range("H4").Value = 0.812868
$ws.Range("I4").Value = 0.05752762181187367
$ws.Range("J4").Value = 0.05752762181187366
$ws.Range("M4").Value = 5.009781333333333
$ws.Range("N4").Value = 15.029344
$ws.Range("O4").Value = 0.3040172587389813
$ws.Range("P4").Value = 0.3040172587389813
$ws.Range("Q4").Value = 1.357430310954667
$ws.Range("R4").Value = 12.216872798592
$ws.Range("S4").Value = 0.01748938988501866
$ws.Range("T4").Value = 0.01748938988501866
$ws.Range("G5").Value = 0.270956
$ws.Range("H5").Value = 0.812868
$ws.Range("I5").Value = 0.05752762181187367
$ws.Range("J5").Value = 0.05752762181187366
$ws.Range("M5").Value = 2.878432333333333
$ws.Range("N5").Value = 8.635297
$ws.Range("O5").Value = 0.1746769068787666
$ws.Range("P5").Value = 0.1746769068787666
$ws.Range("Q5").Value = 0.7799285113106667
$ws.Range("R5").Value = 7.019356601796
$ws.Range("S5").Value = 0.01004874703818956
$ws.Range("T5").Value = 0.01004874703818955
$ws.Range("I6").Value = 0.7443710555244437
$ws.Range("J6").Value = 0.7443710555244437
$ws.Range("M6").Value = 0.8077876666666667
$ws.Range("N6").Value = 2.423363
$ws.Range("O6").Value = 0.04902038147436601
$ws.Range("P6").Value = 0.04902038147436601
$ws.Range("Q6").Value = 2.832103020808222
$ws.Range("R6").Value = 25.488927187274
$ws.Range("S6").Value = 0.03648935310028471
$ws.Range("T6").Value = 0.03648935310028471
$ws.Range("I7").Value = 0.7443710555244437
$ws.Range("J7").Value = 0.7443710555244437
$ws.Range("O7").Value = 0.4722854529078861
$ws.Range("P7").Value = 0.4722854529078861
$ws.Range("S7").Value = 0.3515556210898831
$ws.Range("T7").Value = 0.3515556210898831
$ws.Range("I8").Value = 0.7443710555244437
$ws.Range("J8").Value = 0.7443710555244437
$ws.Range("M8").Value = 5.009781333333333
$ws.Range("N8").Value = 15.029344
$ws.Range("O8").Value = 0.3040172587389813
$ws.Range("P8").Value = 0.3040172587389813
$ws.Range("Q8").Value = 17.56429001481244
$ws.Range("R8").Value = 158.078610133312
$ws.Range("S8").Value = 0.2263016477851834
$ws.Range("T8").Value = 0.2263016477851834
$ws.Range("I9").Value = 0.7443710555244437
$ws.Range("J9").Value = 0.7443710555244437
$ws.Range("M9").Value = 2.878432333333333
$ws.Range("N9").Value = 8.635297
$ws.Range("O9").Value = 0.1746769068787666
$ws.Range("P9").Value = 0.1746769068787666
$ws.Range("Q9").Value = 10.09178184171178
$ws.Range("R9").Value = 90.82603657540599
$ws.Range("S9").Value = 0.1300244335490924
$ws.Range("T9").Value = 0.1300244335490924
$ws.Range("G10").Value = 0.3861913333333333
$ws.Range("H10").Value = 1.158574
$ws.Range("I10").Value = 0.08199364092702592
$ws.Range("J10").Value = 0.08199364092702592
$ws.Range("M10").Value = 0.8077876666666667
$ws.Range("N10").Value = 2.423363
$ws.Range("O10").Value = 0.04902038147436601
$ws.Range("P10").Value = 0.04902038147436601
$ws.Range("Q10").Value = 0.3119605960402222
$ws.Range("R10").Value = 2.807645364362
$ws.Range("S10").Value = 0.004019359556715
$ws.Range("T10").Value = 0.004019359556715001
$ws.Range("G11").Value = 0.3861913333333333
$ws.Range("H11").Value = 1.158574
$ws.Range("I11").Value = 0.08199364092702592
$ws.Range("J11").Value = 0.08199364092702592
$ws.Range("O11").Value = 0.4722854529078861
$ws.Range("P11").Value = 0.4722854529078861
$ws.Range("Q11").Value = 3.005575374139333
$ws.Range("R11").Value = 27.050178367254
$ws.Range("S11").Value = 0.03872440384078702
$ws.Range("T11").Value = 0.03872440384078702
$ws.Range("G12").Value = 0.3861913333333333
$ws.Range("H12").Value = 1.158574
$ws.Range("I12").Value = 0.08199364092702592
$ws.Range("J12").Value = 0.08199364092702592
$ws.Range("M12").Value = 5.009781333333333
$ws.Range("N12").Value = 15.029344
$ws.Range("O12").Value = 0.3040172587389813
$ws.Range("P12").Value = 0.3040172587389813
$ws.Range("Q12").Value = 1.934734132828444
$ws.Range("R12").Value = 17.412607195456
$ws.Range("S12").Value = 0.02492748194866276
$ws.Range("T12").Value = 0.02492748194866276
$ws.Range("G13").Value = 0.3861913333333333
$ws.Range("H13").Value = 1.158574
$ws.Range("I13").Value = 0.08199364092702592
$ws.Range("J13").Value = 0.08199364092702592
$ws.Range("M13").Value = 2.878432333333333
$ws.Range("N13").Value = 8.635297
$ws.Range("O13").Value = 0.1746769068787666
$ws.Range("P13").Value = 0.1746769068787666
$ws.Range("Q13").Value = 1.111625620719778
$ws.Range("R13").Value = 10.004630586478
$ws.Range("S13").Value = 0.01432239558086113
$ws.Range("T13").Value = 0.01432239558086113
$ws.Range("G14").Value = 0.546869
$ws.Range("H14").Value = 1.640607
$ws.Range("I14").Value = 0.1161076817366566
$ws.Range("J14").Value = 0.1161076817366566
$ws.Range("M14").Value = 0.8077876666666667
$ws.Range("N14").Value = 2.423363
$ws.Range("O14").Value = 0.04902038147436601
$ws.Range("P14").Value = 0.04902038147436601
$ws.Range("Q14").Value = 0.4417540334823334
$ws.Range("R14").Value = 3.975786301341
$ws.Range("S14").Value = 0.005691642850835187
$ws.Range("T14").Value = 0.005691642850835188
$ws.Range("G15").Value = 0.546869
$ws.Range("H15").Value = 1.640607
$ws.Range("I15").Value = 0.1161076817366566
$ws.Range("J15").Value = 0.1161076817366566
$ws.Range("O15").Value = 0.4722854529078861
$ws.Range("P15").Value = 0.4722854529078861
$ws.Range("Q15").Value = 4.256066507483
$ws.Range("R15").Value = 38.304598567347
$ws.Range("S15").Value = 0.05483596905508158
$ws.Range("T15").Value = 0.05483596905508158
$ws.Range("G16").Value = 0.546869
$ws.Range("H16").Value = 1.640607
$ws.Range("I16").Value = 0.1161076817366566
$ws.Range("J16").Value = 0.1161076817366566
$ws.Range("M16").Value = 5.009781333333333
$ws.Range("N16").Value = 15.029344
$ws.Range("O16").Value = 0.3040172587389813
$ws.Range("P16").Value = 0.3040172587389813
$ws.Range("Q16").Value = 2.739694107978667
$ws.Range("R16").Value = 24.657246971808
$ws.Range("S16").Value = 0.03529873912011643
$ws.Range("T16").Value = 0.03529873912011643
$ws.Range("G17").Value = 0.546869
$ws.Range("H17").Value = 1.640607
$ws.Range("I17").Value = 0.1161076817366566
$ws.Range("J17").Value = 0.1161076817366566
$ws.Range("M17").Value = 2.878432333333333
$ws.Range("N17").Value = 8.635297
$ws.Range("O17").Value = 0.1746769068787666
$ws.Range("P17").Value = 0.1746769068787666
$ws.Range("Q17").Value = 1.574125411697667
$ws.Range("R17").Value = 14.167128705279
$ws.Range("S17").Value = 0.02028133071062344
$ws.Range("T17").Value = 0.02028133071062344
